$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 436.2857
$ws.Range("J4").Value = 470
$ws.Range("L4").Value = 470
$ws.Range("N4").Value = -698
$ws.Range("H33").Value = 308.92856
$ws.Range("I33").Value = 221.11111
$ws.Range("J33").Value = 467
$ws.Range("K33").Value = 221.11111
$ws.Range("L33").Value = 467
$ws.Range("M33").Value = 7.888890000000004
$ws.Range("N33").Value = -925
$ws.Range("H116").Value = 32954.5
$ws.Range("I116").Value = 32954.5
$ws.Range("K116").Value = 32954.5
$ws.Range("M116").Value = -29512.5
$ws.Range("H132").Value = 20656
$ws.Range("I132").Value = 21063.111
$ws.Range("J132").Value = 6000
$ws.Range("K132").Value = 63189.333
$ws.Range("L132").Value = 18000
$ws.Range("M132").Value = -60659.333
$ws.Range("N132").Value = -23060

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2809.15
$ws.Range("I45").Value = 1741.5
$ws.Range("K45").Value = 1741.5
$ws.Range("M45").Value = -1364.5
$ws.Range("H61").Value = 5750.905
$ws.Range("I61").Value = 935.2105
$ws.Range("K61").Value = 935.2105
$ws.Range("M61").Value = -723.2105
$ws.Range("H74").Value = 359499.47
$ws.Range("I74").Value = 375718.2
$ws.Range("K74").Value = 375718.2
$ws.Range("M74").Value = -374844.2
$ws.Range("H77").Value = 359499.47
$ws.Range("I77").Value = 375718.2
$ws.Range("K77").Value = 1878591
$ws.Range("M77").Value = -1874223
$ws.Range("H132").Value = 1744.4688
$ws.Range("I132").Value = 1358.5769
$ws.Range("K132").Value = 4075.7307
$ws.Range("M132").Value = -1545.7307
$ws.Range("H136").Value = 5750.905
$ws.Range("I136").Value = 935.2105
$ws.Range("K136").Value = 2805.6315
$ws.Range("M136").Value = -255.6315

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 6574.041
$ws.Range("I134").Value = 6662.136
$ws.Range("J134").Value = 5798.8
$ws.Range("K134").Value = 19986.408
$ws.Range("L134").Value = 17396.4
$ws.Range("M134").Value = -17451.408
$ws.Range("N134").Value = -22466.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2767.4546
$ws.Range("I16").Value = 1774.2858
$ws.Range("J16").Value = 4505.5
$ws.Range("K16").Value = 1774.2858
$ws.Range("L16").Value = 4505.5
$ws.Range("M16").Value = -1487.2858
$ws.Range("N16").Value = -5079.5
$ws.Range("H99").Value = 7198.6924
$ws.Range("I99").Value = 4957.4
$ws.Range("K99").Value = 4957.4
$ws.Range("M99").Value = -3459.4
$ws.Range("H113").Value = 2767.4546
$ws.Range("I113").Value = 1774.2858
$ws.Range("J113").Value = 4505.5
$ws.Range("K113").Value = 1774.2858
$ws.Range("L113").Value = 4505.5
$ws.Range("M113").Value = 395.7141999999999
$ws.Range("N113").Value = -8845.5
$ws.Range("H122").Value = 3798.3333
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").Value = ""
$ws.Range("H126").Value = 7198.6924
$ws.Range("I126").Value = 4957.4
$ws.Range("K126").Value = 14872.2
$ws.Range("M126").Value = -12402.2
$ws.Range("H132").Value = 334040.66
$ws.Range("I132").Value = 500505.5
$ws.Range("K132").Value = 1501516.5
$ws.Range("M132").Value = -1498986.5
$ws.Range("H134").Value = 2141.4546
$ws.Range("I134").Value = 1742.7368
$ws.Range("J134").Value = 4666.6665
$ws.Range("K134").Value = 5228.2104
$ws.Range("L134").Value = 13999.9995
$ws.Range("M134").Value = -2693.2104
$ws.Range("N134").Value = -19069.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 185.92857
$ws.Range("I2").Value = 189.33333
$ws.Range("K2").Value = 1135.99998
$ws.Range("M2").Value = -1022.99998
$ws.Range("H26").Value = 152.5
$ws.Range("I26").Value = 150
$ws.Range("J26").Value = 170
$ws.Range("K26").Value = 450
$ws.Range("L26").Value = 510
$ws.Range("M26").Value = -162
$ws.Range("N26").Value = -1086
$ws.Range("H37").Value = 42089.176
$ws.Range("J37").Value = 42089.176
$ws.Range("L37").Value = 126267.528
$ws.Range("N37").Value = -126491.528
$ws.Range("H40").Value = 30.6
$ws.Range("I40").Value = 30.6
$ws.Range("K40").Value = 122.4
$ws.Range("M40").Value = -53.40000000000001
$ws.Range("H86").Value = 1527.7142
$ws.Range("I86").Value = 1498.3334
$ws.Range("J86").Value = 1549.75
$ws.Range("K86").Value = 4495.0002
$ws.Range("L86").Value = 4649.25
$ws.Range("M86").Value = -3309.0002
$ws.Range("N86").Value = -7021.25
$ws.Range("H89").Value = 1527.7142
$ws.Range("I89").Value = 1498.3334
$ws.Range("J89").Value = 1549.75
$ws.Range("K89").Value = 13485.0006
$ws.Range("L89").Value = 13947.75
$ws.Range("M89").Value = -7557.000599999999
$ws.Range("N89").Value = -25803.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 37999.8
$ws.Range("I52").Value = 33333
$ws.Range("K52").Value = 33333
$ws.Range("M52").Value = -33074
$ws.Range("H132").Value = 2827.2632
$ws.Range("I132").Value = 1979.0769
$ws.Range("K132").Value = 5937.2307
$ws.Range("M132").Value = -3407.2307

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3852.1428
$ws.Range("I7").Value = 3852.1428
$ws.Range("K7").Value = 3852.1428
$ws.Range("M7").Value = -3740.1428
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").Value = ""
$ws.Range("H18").Value = 19996.666
$ws.Range("J18").Value = 19996.666
$ws.Range("L18").Value = 19996.666
$ws.Range("N18").Value = -20340.666
$ws.Range("H22").Value = 1776
$ws.Range("I22").Value = 1494.375
$ws.Range("J22").Value = 2057.625
$ws.Range("K22").Value = 1494.375
$ws.Range("L22").Value = 2057.625
$ws.Range("M22").Value = -1199.375
$ws.Range("N22").Value = -2647.625
$ws.Range("H27").Value = 1776
$ws.Range("I27").Value = 1494.375
$ws.Range("J27").Value = 2057.625
$ws.Range("K27").Value = 1494.375
$ws.Range("L27").Value = 2057.625
$ws.Range("M27").Value = -1387.375
$ws.Range("N27").Value = -2271.625
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").Value = ""
$ws.Range("H40").Value = 2566.5557
$ws.Range("I40").Value = 2262.375
$ws.Range("K40").Value = 2262.375
$ws.Range("M40").Value = -2126.375
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("M48").Value = ""
$ws.Range("N48").Value = ""
$ws.Range("H122").Value = 5849.1665
$ws.Range("I122").Value = 4033.3333
$ws.Range("K122").Value = 12099.9999
$ws.Range("M122").Value = -9649.999899999999
$ws.Range("H125").Value = 86715
$ws.Range("J125").Value = 86715
$ws.Range("L125").Value = 86715
$ws.Range("N125").Value = -96555
$ws.Range("H126").Value = 3852.1428
$ws.Range("I126").Value = 3852.1428
$ws.Range("K126").Value = 11556.4284
$ws.Range("M126").Value = -9086.428400000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H115").Value = 39772.727
$ws.Range("J115").Value = 39772.727
$ws.Range("L115").Value = 39772.727
$ws.Range("N115").Value = -42906.727
$ws.Range("H122").Value = 44524.234
$ws.Range("I122").Value = 55068.555
$ws.Range("K122").Value = 165205.665
$ws.Range("M122").Value = -162755.665
